# Update "Estado de Cuenta" worksheet data (rows 16-29) with the new
# database ordering: interleave LIBARDO / YOJHANA rows by ascending
# period (1904 -> 1910) instead of grouping all of LIBARDO's periods
# followed by all of YOJHANA's periods. Also update YOJHANA's
# "Salario Basico" (column G) to the new value 781242.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$libardoDoc    = "73091597"
$libardoName   = "LIBARDO CRISTOBAL CASTELLAR CASTELLAR"
$yojhanaDoc    = "50926657"
$yojhanaName   = "YOJHANA LUCIA ALVAREZ CABRALES"

# row => (DocType, DocNumber, Name, Period, ValorMora, SalarioBasico)
$rows = @(
    @{ Row=16; Tipo="CC"; Doc=$libardoDoc; Nombre=$libardoName; Periodo="1904"; Mora=33125; Salario=828116 },
    @{ Row=17; Tipo="CC"; Doc=$yojhanaDoc; Nombre=$yojhanaName; Periodo="1904"; Mora=33125; Salario=781242 },
    @{ Row=18; Tipo="CC"; Doc=$libardoDoc; Nombre=$libardoName; Periodo="1905"; Mora=33125; Salario=828116 },
    @{ Row=19; Tipo="CC"; Doc=$yojhanaDoc; Nombre=$yojhanaName; Periodo="1905"; Mora=33125; Salario=781242 },
    @{ Row=20; Tipo="CC"; Doc=$libardoDoc; Nombre=$libardoName; Periodo="1906"; Mora=33125; Salario=828116 },
    @{ Row=21; Tipo="CC"; Doc=$yojhanaDoc; Nombre=$yojhanaName; Periodo="1906"; Mora=33125; Salario=781242 },
    @{ Row=22; Tipo="CC"; Doc=$libardoDoc; Nombre=$libardoName; Periodo="1907"; Mora=33125; Salario=828116 },
    @{ Row=23; Tipo="CC"; Doc=$yojhanaDoc; Nombre=$yojhanaName; Periodo="1907"; Mora=33125; Salario=781242 },
    @{ Row=24; Tipo="CC"; Doc=$libardoDoc; Nombre=$libardoName; Periodo="1908"; Mora=33125; Salario=828116 },
    @{ Row=25; Tipo="CC"; Doc=$yojhanaDoc; Nombre=$yojhanaName; Periodo="1908"; Mora=33125; Salario=781242 },
    @{ Row=26; Tipo="CC"; Doc=$libardoDoc; Nombre=$libardoName; Periodo="1909"; Mora=33125; Salario=828116 },
    @{ Row=27; Tipo="CC"; Doc=$yojhanaDoc; Nombre=$yojhanaName; Periodo="1909"; Mora=31249; Salario=781242 },
    @{ Row=28; Tipo="CC"; Doc=$libardoDoc; Nombre=$libardoName; Periodo="1910"; Mora=23187; Salario=828116 },
    @{ Row=29; Tipo="CC"; Doc=$yojhanaDoc; Nombre=$yojhanaName; Periodo="1910"; Mora=21874; Salario=781242 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value2 = $item.Tipo
    $ws.Cells.Item($r, 3).Value2 = $item.Doc
    $ws.Cells.Item($r, 4).Value2 = $item.Nombre
    $ws.Cells.Item($r, 5).Value2 = $item.Periodo
    $ws.Cells.Item($r, 6).Value2 = $item.Mora
    $ws.Cells.Item($r, 7).Value2 = $item.Salario
}
